$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (row 3)
$ws.Range("A3").Value = "Kaidi"
$ws.Range("B3").Value = "Rim"
$ws.Range("C3").Value = 22010892
$ws.Range("D3").Value = "kaidirim12498@gmail.com"
$ws.Range("E3").Value = "super !"

# Turn the email into a mailto hyperlink, matching the existing D2 cell
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:kaidirim12498@gmail.com")
$ws.Range("D3").Style = "Lien hypertexte"
